# feat(command): add a command gathercolumn corresponding to the new function.
# Gathers (duplicates) the 12-column block I:T (one full "Alain/Henri/Tony/
# Dulcinee" x OUI/NON cycle) and inserts a fresh copy of it just before the
# existing "Adresse de courriel" / trailing empty column pair that currently
# sit at HY:HZ, shifting that pair two columns further right (-> IK:IL).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$ws.Columns("I:T").Copy()
$ws.Columns("HY:IJ").Insert($xlShiftToRight)
